$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values that changed
$ws.Range("B10").Value = 0.7877551020408163
$ws.Range("B11").Value = 0.1619047619047619
$ws.Range("B12").Value = 0.01730103806228374
$ws.Range("B17").Value = 0.4443907156673114
$ws.Range("B19").Value = 0.08380952380952381
$ws.Range("B20").Value = 0.4419405931287398

# Add new rows 21-27
$ws.Range("A21").Value = "Flex Monte Carlo"
$ws.Range("B21").Value = 1

$ws.Range("A22").Value = "Flex Netzreserve"
$ws.Range("B22").Value = 1

$ws.Range("A23").Value = "Flex Erfolgreiche OPP"
$ws.Range("B23").Value = 0

$ws.Range("A24").Value = "Flex Reserve krit Leitungen"
$ws.Range("B24").Value = 1

$ws.Range("A25").Value = "Flexibilität Gesamt"
$ws.Range("B25").Value = 0.75

$ws.Range("A26").Value = "Buffer Capacity"
$ws.Range("B26").Value = 0.03024019443551216

$ws.Range("A27").Value = "Feasible operating region"
$ws.Range("B27").Value = 0.09318236976738642
